# Translate "disciplina" course catalogue sheet from Portuguese to Spanish,
# keep only the first few course entries (CIC, FM, ING, FPOO) and clear the
# remaining rows of data, per the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: CIC / Conceitos de internet das Coisas -> Conceptos de internet de las Cosas
$ws.Range("B1").Value = "Conceptos de internet de las Cosas"

# Row 2: FM / Fundamentos de Matemática -- unchanged

# Row 3: ING / Inglês -> Ingles
$ws.Range("B3").Value = "Ingles"

# Row 4: used to hold IPI / Introdução à Programação para Internet.
# It's replaced with the FPOO entry (translated), and the old IPI / PWIOT /
# SGBD rows that used to sit between them are dropped entirely.
$ws.Range("A4").Value = "FPOO"
$ws.Range("B4").Value = "Fundamentos de Programación Orientada a Objetos"

# Rows 5-49 held the rest of the course catalogue (PM, PCMR, SRSI, ... EST);
# all of that detail is removed, leaving only blank, styled cells behind.
$ws.Range("A5:B49").ClearContents()
